$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" '27.177.47'
Set-TextValue "E2" '  -0.11%  '
Set-TextValue "D3" '1.900.52'
Set-TextValue "E3" '  -0.24%  '
Set-TextValue "D4" '1.003'
Set-TextValue "E4" '  +0.29%  '
Set-TextValue "D5" '307.32'
Set-TextValue "E6" '  +0.25%  '
Set-TextValue "D7" '0.5246'
Set-TextValue "E7" '  -0.09%  '
Set-TextValue "D8" '0.3809'
Set-TextValue "E8" '  +0.80%  '
Set-TextValue "D9" '0.07293'
Set-TextValue "E9" '  +0.45%  '
Set-TextValue "D10" '21.38'
Set-TextValue "E10" '  +1.12%  '
Set-TextValue "D11" '0.9041'
Set-TextValue "E11" '  +0.48%  '
Set-TextValue "D12" '0.08164'
Set-TextValue "E12" '  -3.18%  '
Set-TextValue "D13" '95.47'
Set-TextValue "E13" '  +0.77%  '
Set-TextValue "D14" '5.349'
Set-TextValue "E14" '  +1.49%  '
Set-TextValue "D15" '1.799.02'
Set-TextValue "E15" '  -5.46%  '
Set-TextValue "E16" '  +0.28%  '
Set-TextValue "D17" '0.000008654'
Set-TextValue "E17" '  +0.28%  '
Set-TextValue "D18" '14.69'
Set-TextValue "E18" '  +0.85%  '
Set-TextValue "D19" '1.002'
Set-TextValue "E19" '  +0.21%  '
Set-TextValue "D20" '27.204.29'
Set-TextValue "E20" '  -0.11%  '
Set-TextValue "D21" '5.098'
Set-TextValue "E21" '  +0.74%  '
Set-TextValue "D22" '10.79'
Set-TextValue "E22" '  +1.92%  '
Set-TextValue "D23" '6.453'
Set-TextValue "E23" '  +0.28%  '
Set-TextValue "B24" 'LidoDAOToken'
Set-TextValue "C24" 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue "D24" '2.330'
Set-TextValue "E24" '  +2.50%  '
Set-TextValue "B25" 'Monero'
Set-TextValue "C25" 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue "D25" '149.49'
Set-TextValue "E25" '  +1.78%  '
Set-TextValue "E26" '  +0.25%  '
Set-TextValue "E27" '  -0.44%  '
Set-TextValue "D28" '115.96'
Set-TextValue "E28" '  +0.99%  '
Set-TextValue "D29" '4.830'
Set-TextValue "E29" '  +0.49%  '
Set-TextValue "D30" '4.881'
Set-TextValue "E30" '  -0.98%  '
Set-TextValue "D31" '0.09234'
Set-TextValue "E31" '  -0.67%  '
Set-TextValue "D32" '0.05058'
Set-TextValue "E32" '  -0.15%  '
Set-TextValue "D33" '0.7927'
Set-TextValue "E33" '  -1.87%  '
Set-TextValue "D34" '1.225'
Set-TextValue "E34" '  -0.91%  '
Set-TextValue "D35" '2.981'
Set-TextValue "E35" '  +0.98%  '
Set-TextValue "D36" '3.367'
Set-TextValue "E36" '  +0.33%  '
Set-TextValue "D37" '2.649'
Set-TextValue "E37" '  +1.74%  '
Set-TextValue "D38" '0.5708'
Set-TextValue "E38" '  -0.22%  '
Set-TextValue "D39" '0.01988'
Set-TextValue "E39" '  +0.02%  '
Set-TextValue "D40" '1.081'
Set-TextValue "D41" '9.014'
Set-TextValue "E41" '  +0.53%  '
Set-TextValue "D42" '6.587'
Set-TextValue "E42" '  -0.94%  '
Set-TextValue "D43" '116.29'
Set-TextValue "E43" '  -1.19%  '
Set-TextValue "D44" '0.1512'
Set-TextValue "E44" '  -0.19%  '
Set-TextValue "D45" '0.4875'
Set-TextValue "E45" '  +0.55%  '
Set-TextValue "E46" '  +0.30%  '
Set-TextValue "D47" '10.10'
Set-TextValue "E47" '  -0.48%  '
Set-TextValue "D48" '1.628'
Set-TextValue "E48" '  +0.81%  '
Set-TextValue "D49" '38.51'
Set-TextValue "E49" '  +2.83%  '
Set-TextValue "D50" '63.91'
Set-TextValue "E50" '  +0.47%  '
Set-TextValue "D51" '0.05960'
Set-TextValue "E51" '  +0.42%  '

Write-Output "Updated cryptos list"
